$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.506.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.914.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'244.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.68%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.01%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4810"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.40%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.44%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06723"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.21%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'110.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.80%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'19.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.43%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.912.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.94%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07560"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.63%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6674"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.32%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'293.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.501.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.24%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007580"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.21%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.164.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.47%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.465"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.22%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.396"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.19%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.456"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.49%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'164.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.32%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'20.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -6.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.125"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.90%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.1068"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.403"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.83%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.159"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.021"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.94%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04967"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.78%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7288"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.35%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.77%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'VeChain"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.02057"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'HuobiToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'2.742"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.08%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Frax"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.9996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.10%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.669"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.39%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'110.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.008"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4410"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.8623"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.97%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'5.898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.31%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.9998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'68.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.38%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'49.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.18%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.1236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.39%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.2520"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.62%  "
$ws.Range("E51").Style = "Normal"
